$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 590.1667
$ws.Range("I19").Value = 357
$ws.Range("K19").Value = 357
$ws.Range("M19").Value = -182
$ws.Range("H40").Value = 1000
$ws.Range("I40").Value = 1000
$ws.Range("J40").Value = 1000
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 1000
$ws.Range("M40").Value = -825
$ws.Range("N40").Value = -1350
$ws.Range("H74").Value = 3966.7144
$ws.Range("I74").Value = 3761.125
$ws.Range("J74").Value = 4093.2307
$ws.Range("K74").Value = 3761.125
$ws.Range("L74").Value = 4093.2307
$ws.Range("M74").Value = -2825.125
$ws.Range("N74").Value = -5965.2307
$ws.Range("H76").Value = 3473
$ws.Range("I76").Value = 3422.5557
$ws.Range("K76").Value = 3422.5557
$ws.Range("M76").Value = -3107.5557
$ws.Range("H77").Value = 3966.7144
$ws.Range("I77").Value = 3761.125
$ws.Range("J77").Value = 4093.2307
$ws.Range("K77").Value = 18805.625
$ws.Range("L77").Value = 20466.1535
$ws.Range("M77").Value = -14125.625
$ws.Range("N77").Value = -29826.1535
$ws.Range("H79").Value = 3473
$ws.Range("I79").Value = 3422.5557
$ws.Range("K79").Value = 3422.5557
$ws.Range("M79").Value = -2330.5557
$ws.Range("H135").Value = 129990.75
$ws.Range("I135").Value = 171987.67
$ws.Range("J135").Value = 4000
$ws.Range("K135").Value = 1547889.03
$ws.Range("L135").Value = 36000
$ws.Range("M135").Value = -1545354.03
$ws.Range("N135").Value = -41070
$ws.Range("H138").Value = 3127478.5
$ws.Range("I138").Value = 2551.6667
$ws.Range("J138").Value = 4084088.8
$ws.Range("K138").Value = 7655.000100000001
$ws.Range("L138").Value = 12252266.4
$ws.Range("M138").Value = -2515.000100000001
$ws.Range("N138").Value = -12262546.4
$ws.Range("H141").Value = 2238.5881
$ws.Range("I141").Value = 1531.9259
$ws.Range("J141").Value = 4964.2856
$ws.Range("K141").Value = 4595.7777
$ws.Range("L141").Value = 14892.8568
$ws.Range("M141").Value = 584.2223000000004
$ws.Range("N141").Value = -25252.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23371.445
$ws.Range("I32").Value = 23370.31
$ws.Range("J32").Value = 23374.322
$ws.Range("K32").Value = 23370.31
$ws.Range("L32").Value = 23374.322
$ws.Range("M32").Value = -23083.31
$ws.Range("N32").Value = -23948.322
$ws.Range("H63").Value = 4499.3335
$ws.Range("I63").Value = 4856.2856
$ws.Range("J63").Value = 3250
$ws.Range("K63").Value = 4856.2856
$ws.Range("L63").Value = 3250
$ws.Range("M63").Value = -4170.2856
$ws.Range("N63").Value = -4622
$ws.Range("H66").Value = 4499.3335
$ws.Range("I66").Value = 4856.2856
$ws.Range("J66").Value = 3250
$ws.Range("K66").Value = 24281.428
$ws.Range("L66").Value = 16250
$ws.Range("M66").Value = -20849.428
$ws.Range("N66").Value = -23114
$ws.Range("H88").Value = 5143.3184
$ws.Range("I88").Value = 2253
$ws.Range("J88").Value = 5785.6113
$ws.Range("K88").Value = 2253
$ws.Range("L88").Value = 5785.6113
$ws.Range("M88").Value = -1847
$ws.Range("N88").Value = -6597.6113
$ws.Range("H91").Value = 5143.3184
$ws.Range("I91").Value = 2253
$ws.Range("J91").Value = 5785.6113
$ws.Range("K91").Value = 2253
$ws.Range("L91").Value = 5785.6113
$ws.Range("M91").Value = -849
$ws.Range("N91").Value = -8593.6113

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1370.8
$ws.Range("I20").Value = 968
$ws.Range("J20").Value = 1723.25
$ws.Range("K20").Value = 968
$ws.Range("L20").Value = 1723.25
$ws.Range("M20").Value = -721
$ws.Range("N20").Value = -2217.25
$ws.Range("H86").Value = 15763.5625
$ws.Range("I86").Value = 18707.691
$ws.Range("K86").Value = 18707.691
$ws.Range("M86").Value = -17584.691
$ws.Range("H89").Value = 15763.5625
$ws.Range("I89").Value = 18707.691
$ws.Range("K89").Value = 93538.45499999999
$ws.Range("M89").Value = -87922.45499999999
$ws.Range("H105").Value = 71433256
$ws.Range("I105").Value = 125003950
$ws.Range("J105").Value = 5666.6665
$ws.Range("K105").Value = 125003950
$ws.Range("L105").Value = 5666.6665
$ws.Range("M105").Value = -125002203
$ws.Range("N105").Value = -9160.666499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 252952.94
$ws.Range("I31").Value = 73661.21000000001
$ws.Range("J31").Value = 326778.94
$ws.Range("K31").Value = 73661.21000000001
$ws.Range("L31").Value = 326778.94
$ws.Range("M31").Value = -73366.21000000001
$ws.Range("N31").Value = -327368.94
$ws.Range("H34").Value = 252952.94
$ws.Range("I34").Value = 73661.21000000001
$ws.Range("J34").Value = 326778.94
$ws.Range("K34").Value = 73661.21000000001
$ws.Range("L34").Value = 326778.94
$ws.Range("M34").Value = -73459.21000000001
$ws.Range("N34").Value = -327182.94
$ws.Range("H62").Value = 3082.85
$ws.Range("I62").Value = 2876.125
$ws.Range("K62").Value = 2876.125
$ws.Range("M62").Value = -2252.125
$ws.Range("H65").Value = 3082.85
$ws.Range("I65").Value = 2876.125
$ws.Range("K65").Value = 14380.625
$ws.Range("M65").Value = -11260.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 39506.62
$ws.Range("I70").Value = 46516.332
$ws.Range("J70").Value = 5860
$ws.Range("K70").Value = 46516.332
$ws.Range("L70").Value = 5860
$ws.Range("M70").Value = -46246.332
$ws.Range("N70").Value = -6400
$ws.Range("H73").Value = 39506.62
$ws.Range("I73").Value = 46516.332
$ws.Range("J73").Value = 5860
$ws.Range("K73").Value = 46516.332
$ws.Range("L73").Value = 5860
$ws.Range("M73").Value = -45580.332
$ws.Range("N73").Value = -7732
$ws.Range("H80").Value = 4144
$ws.Range("I80").Value = 2600
$ws.Range("K80").Value = 2600
$ws.Range("M80").Value = -1602
$ws.Range("H83").Value = 4144
$ws.Range("I83").Value = 2600
$ws.Range("K83").Value = 13000
$ws.Range("M83").Value = -8008
$ws.Range("H97").Value = 1764.8966
$ws.Range("I97").Value = 1713
$ws.Range("J97").Value = 1880.2222
$ws.Range("K97").Value = 1713
$ws.Range("L97").Value = 1880.2222
$ws.Range("M97").Value = -1217
$ws.Range("N97").Value = -2872.2222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 74454.13
$ws.Range("I132").Value = 6780.7
$ws.Range("J132").Value = 209801
$ws.Range("K132").Value = 20342.1
$ws.Range("L132").Value = 629403
$ws.Range("M132").Value = -17812.1
$ws.Range("N132").Value = -634463
$ws.Range("H136").Value = 58363.54
$ws.Range("I136").Value = 41240.777
$ws.Range("J136").Value = 104595
$ws.Range("K136").Value = 123722.331
$ws.Range("L136").Value = 313785
$ws.Range("M136").Value = -121172.331
$ws.Range("N136").Value = -318885

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2048
$ws.Range("J81").Value = 2386.5715
$ws.Range("L81").Value = 4773.143
$ws.Range("N81").Value = -6895.143
$ws.Range("H84").Value = 2048
$ws.Range("J84").Value = 2386.5715
$ws.Range("L84").Value = 23865.715
$ws.Range("N84").Value = -34473.715
$ws.Range("H132").Value = 46066.38
$ws.Range("I132").Value = 46654.137
$ws.Range("J132").Value = 45504.176
$ws.Range("K132").Value = 139962.411
$ws.Range("L132").Value = 136512.528
$ws.Range("M132").Value = -137432.411
$ws.Range("N132").Value = -141572.528
